$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("generator_file")

# Insert a new column before column C (old C/D become D/E).
$ws.Columns.Item(3).Insert()

# The insert operation copies formatting from column B into the new,
# still-empty column C cells. Clear those placeholder cells in the data
# rows that shouldn't carry a value/format in column C.
$ws.Range("C2:C6").Clear()

# New column C width/format matches column B (no bestFit, explicit width).
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth()

# Header row: C1 = "prime_mover" (bold + text number format, like B1).
$ws.Range("C1").Value = "prime_mover"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").NumberFormat = "@"

# Row 7 (plant 55350, generator_id 1): was "update prime_mover to CT",
# now specifies prime_mover CA directly and swaps generator_id to 3.
$ws.Range("C7").Value = "CA"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("D7").Value = "generator_id"
$ws.Range("E7").Value = "3"
$ws.Range("E7").NumberFormat = "@"

# Row 8 (plant 55350, generator_id 3): was "update prime_mover to CA",
# now specifies prime_mover CT directly and swaps generator_id to 1.
$ws.Range("C8").Value = "CT"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("D8").Value = "generator_id"
$ws.Range("E8").Value = "1"
$ws.Range("E8").NumberFormat = "@"

# Restore the active selection shown in the saved file.
$ws.Range("C2").Select() | Out-Null
